$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Status column (D) updates -------------------------------------------------
$ws.Range("D86").Value = "S"
$ws.Range("D89").Value = "S"
$ws.Range("D93").Value = "PS"
$ws.Range("D97").Value = "PS"
$ws.Range("D98").Value = "PS"
$ws.Range("D99").Value = "PS"
$ws.Range("D100").Value = "PS"
$ws.Range("D101").Value = "PS"
$ws.Range("D103").Value = "V"
$ws.Range("D105").Value = "PS"

# --- New comment in "Modalites de verification" column (G) ---------------------
$ws.Range("G106").Value = "Voir le registre des livrables"

# --- Selection / view state ------------------------------------------------------
$ws.Range("D103").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 102
